$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "Miss Dina Nasr, Administrator"
$new = "Administrator, Miss Dina Nasr"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$colG = $ws.Range("G1:G$lastRow")
$result = $colG.Replace($old, $new, 1)
